$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "French" column header (match the style of the other header cells)
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "French"

# Add French translations for existing rows
$ws.Range("G2").Value = "(étiquette)"
$ws.Range("G3").Value = "(mod``ele)"
$ws.Range("G4").Value = "surface"

# Add new "freehand" row with Polish + French translations
$ws.Range("A5").Value = "freehand"
$ws.Range("F5").Value = "Pismo odręczne"
$ws.Range("G5").Value = "(écriture)"

# Add new license-tier rows
$ws.Range("A6").Value = "Basic"
$ws.Range("A7").Value = "Pro"
$ws.Range("A8").Value = "Business"
$ws.Range("G8").Value = "Enterprise"

# Update the active selection to mirror the new next-empty-row
$ws.Range("A9").Select()
